# The commit swaps the two theme parts in the deck: the theme that used to
# live in ppt/theme/theme1.xml (the "Integral" / Red Violet design applied to
# the Slide Master, and therefore to every slide) ends up with the colours
# that used to live in ppt/theme/theme2.xml (the default "Office Theme" that
# the Notes Master uses), and vice versa.
#
# The only thing that actually differs between the two theme parts is the
# <a:clrScheme> (its name, and ten of its twelve colours - dk1/black and
# lt1/white are shared) - the font scheme and the format scheme are byte for
# byte identical in both files. So the visible effect of the swap is fully
# captured by re-pointing the Slide Master's theme colour scheme at the
# palette that used to belong to the Office Theme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

function ToRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# msoThemeColorDark1 .. msoThemeColorFollowedHyperlink (indices 1-12),
# set to the "Office Theme" palette (previously ppt/theme/theme2.xml).
$colors.Item(1).RGB  = ToRgb 0x00 0x00 0x00   # Dark 1
$colors.Item(2).RGB  = ToRgb 0xFF 0xFF 0xFF   # Light 1
$colors.Item(3).RGB  = ToRgb 0x44 0x54 0x6A   # Dark 2
$colors.Item(4).RGB  = ToRgb 0xE7 0xE6 0xE6   # Light 2
$colors.Item(5).RGB  = ToRgb 0x5B 0x9B 0xD5   # Accent 1
$colors.Item(6).RGB  = ToRgb 0xED 0x7D 0x31   # Accent 2
$colors.Item(7).RGB  = ToRgb 0xA5 0xA5 0xA5   # Accent 3
$colors.Item(8).RGB  = ToRgb 0xFF 0xC0 0x00   # Accent 4
$colors.Item(9).RGB  = ToRgb 0x44 0x72 0xC4   # Accent 5
$colors.Item(10).RGB = ToRgb 0x70 0xAD 0x47   # Accent 6
$colors.Item(11).RGB = ToRgb 0x05 0x63 0xC1   # Hyperlink
$colors.Item(12).RGB = ToRgb 0x95 0x4F 0x72   # Followed Hyperlink
